$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count

    # Swap column A (약품코드/code) and column C (번호/number) contents for every row,
    # including the header row, so the new layout becomes:
    # A=번호, B=약품명, C=약품코드, D=유효기간
    for ($r = 1; $r -le $lastRow; $r++) {
        $a = $ws.Cells.Item($r, 1).Value2
        $c = $ws.Cells.Item($r, 3).Value2
        $ws.Cells.Item($r, 1).Value2 = $c
        $ws.Cells.Item($r, 3).Value2 = $a
    }

    # Column widths: col A now holds short numbers (width 4), col B stays 14,
    # col C now holds the longer codes (width 6, matching col D).
    $ws.Columns.Item(1).ColumnWidth = 3.1666666666666667
    $ws.Columns.Item(3).ColumnWidth = 5.166666666666667
}
